$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.90"
$ws.Range("E2").Value = "'2.83%"
$ws.Range("D3").Value = "'41.21"
$ws.Range("E3").Value = "'2.05%"
$ws.Range("E4").Value = "'-0.55%"
$ws.Range("D5").Value = "'0.07501"
$ws.Range("E5").Value = "'2.71%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.363"
$ws.Range("E6").Value = "'1.97%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.572"
$ws.Range("E7").Value = "'3.11%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9263"
$ws.Range("E8").Value = "'0.84%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.401"
$ws.Range("E9").Value = "'0.17%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1221"
$ws.Range("E10").Value = "'2.72%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1837"
$ws.Range("E11").Value = "'6.73%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08836"
$ws.Range("E12").Value = "'2.58%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04104"
$ws.Range("E13").Value = "'-1.39%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("D15").Value = "'0.005871"
$ws.Range("E15").Value = "'0.95%"
$ws.Range("E16").Value = "'-1.52%"
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3288"
$ws.Range("E17").Value = "'0.08%"
$ws.Range("B18").Value = "MCDex"
$ws.Range("C18").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D18").Value = "'8.012"
$ws.Range("E18").Value = "'2.80%"
$ws.Range("B19").Value = "ProBitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D19").Value = "'0.1418"
$ws.Range("E19").Value = "'5.00%"
$ws.Range("B20").Value = "ZBToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D20").Value = "'0.2964"
$ws.Range("E20").Value = "'2.92%"
$ws.Range("B21").Value = "BitForexToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D21").Value = "'0.001281"
$ws.Range("E21").Value = "'3.09%"
$ws.Range("D22").Value = "'0.04042"
$ws.Range("E22").Value = "'5.00%"
$ws.Range("D23").Value = "'0.001266"
$ws.Range("E23").Value = "'-0.08%"
$ws.Range("D24").Value = "'0.003876"
$ws.Range("E24").Value = "'0.95%"
$ws.Range("D25").Value = "'0.0001230"
$ws.Range("E25").Value = "'-3.86%"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("D38").Value = "'0.02417"
$ws.Range("E38").Value = "'4.15%"
$ws.Range("D39").Value = "'0.05222"
$ws.Range("E39").Value = "'4.98%"
$ws.Range("D40").Value = "'0.005991"
$ws.Range("E40").Value = "'-6.59%"
$ws.Range("D41").Value = "'0.007793"
$ws.Range("E41").Value = "'1.57%"
$ws.Range("D42").Value = "'0.1325"
$ws.Range("E42").Value = "'3.98%"
$ws.Range("D43").Value = "'0.007362"
$ws.Range("E43").Value = "'0.21%"
$ws.Range("D44").Value = "'0.008120"
$ws.Range("E44").Value = "'15.24%"
$ws.Range("D45").Value = "'0.2969"
$ws.Range("E45").Value = "'-4.94%"
$ws.Range("D46").Value = "'0.00006252"
$ws.Range("E46").Value = "'-2.76%"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.04511"
$ws.Range("E48").Value = "'-82.00%"
$ws.Range("D49").Value = "'0.004197"
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.04%"
